$wb = $excel.ActiveWorkbook

# Rename the "ExternalReferences_200" sheet (tab) to "Links_200"
$extRefSheet = $wb.Worksheets.Item("ExternalReferences_200")
$extRefSheet.Name = "Links_200"

# Update the CodeSchemes sheet header/value referencing the renamed sheet
$codeSchemes = $wb.Worksheets.Item("CodeSchemes")
$codeSchemes.Range("R1").Value = "LINKSSHEET"
$codeSchemes.Range("R2").Value = "Links_200"

# Column width adjustments on CodeSchemes sheet (A and B)
# (ColumnWidth values chosen to land on the closest representable stored
# width given the engine's internal rounding: ~14.43 and ~41.29)
$codeSchemes.Columns.Item(1).ColumnWidth = 13.65
$codeSchemes.Columns.Item(2).ColumnWidth = 40.5
